$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 143, pushing the existing rows
# 143..245 down to 144..246 (dimension grows from R245 to R246).
$ws.Rows(143).Insert()

# Populate the newly inserted row 143 with its data.
$ws.Range("A143").Value = 11
$ws.Range("B143").Value = "Vega Monumental Concepción"
$ws.Range("C143").Value = "Bíobío"
$ws.Range("D143").Value = 45216
$ws.Range("E143").Value = 8
$ws.Range("F143").Value = 100112043
$ws.Range("G143").Value = "Pepino ensalada"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 100
$ws.Range("K143").Value = 14000
$ws.Range("L143").Value = 15000
$ws.Range("M143").Value = 14500
$ws.Range("N143").Value = "$/caja 60 unidades"
$ws.Range("O143").Value = "Región de Arica y Parinacota"
$ws.Range("P143").Value = 242
$ws.Range("Q143").Value = 60
$ws.Range("R143").Value = "Hortaliza"
